# Update Work Week and Social Spending
# (Guatemala GDP per Capita data refresh: refreshed values for 1920-2010
#  and six new rows appended for 2011-2016.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Map of Year -> new Data value (kept as text, matching the column's
# original shared-string/text storage rather than becoming a number).
$values = @{
    1920 = "2028"
    1921 = "2203"
    1922 = "2058"
    1923 = "2240"
    1924 = "2397"
    1925 = "2327"
    1926 = "2327"
    1927 = "2453"
    1928 = "2482"
    1929 = "2742"
    1930 = "2831"
    1931 = "2582"
    1932 = "2200"
    1933 = "2165"
    1934 = "2410"
    1935 = "2729"
    1936 = "3674"
    1937 = "3516"
    1938 = "3545"
    1939 = "3916"
    1940 = "4371"
    1941 = "4503"
    1942 = "4463"
    1943 = "2925"
    1944 = "2775"
    1945 = "2761"
    1946 = "3191"
    1947 = "3148"
    1948 = "3169"
    1949 = "3365"
    1950 = "3323"
    1951 = "3263"
    1952 = "3226"
    1953 = "3239"
    1954 = "3194"
    1955 = "3170"
    1956 = "3347"
    1957 = "3424"
    1958 = "3470"
    1959 = "3524"
    1960 = "3496"
    1961 = "3532"
    1962 = "3543"
    1963 = "3762"
    1964 = "3814"
    1965 = "3900"
    1966 = "4033"
    1967 = "4114"
    1968 = "4385"
    1969 = "4497"
    1970 = "4653"
    1971 = "4801"
    1972 = "5037"
    1973 = "5255"
    1974 = "5461"
    1975 = "5439"
    1976 = "5705"
    1977 = "6005"
    1978 = "6156"
    1979 = "6298"
    1980 = "6384"
    1981 = "6256"
    1982 = "5903"
    1983 = "5585"
    1984 = "5469"
    1985 = "5292"
    1986 = "5136"
    1987 = "5098"
    1988 = "5086"
    1989 = "5144"
    1990 = "5165"
    1991 = "5166.54986601267"
    1992 = "5303.59854742543"
    1993 = "5382.9787121035"
    1994 = "5414.54115356924"
    1995 = "5501.04222678237"
    1996 = "5511.90470568081"
    1997 = "5601.20750520503"
    1998 = "5730.64575362661"
    1999 = "5824.88328486933"
    2000 = "5860.2718525717"
    2001 = "5882.71970593542"
    2002 = "5980.24200772415"
    2003 = "6004.14467312046"
    2004 = "6070.19423096165"
    2005 = "6130.48233756792"
    2006 = "6305.81859145986"
    2007 = "6545.75235326108"
    2008 = "6604.68242304307"
    2009 = "6488.64623816622"
    2010 = "6526.45756959102"
    2011 = "6650"
    2012 = "6714"
    2013 = "6829"
    2014 = "6981"
    2015 = "7138"
    2016 = "7221"
}

# Format the whole Data column (existing rows plus the new ones we are
# about to append) as Text first so the numeric-looking values below are
# stored as text, matching the original column's cell type.
$ws.Range("E2:E98").NumberFormat = "@"

# Refresh the existing rows (2-92, years 1920-2010) in place.
for ($r = 2; $r -le 92; $r++) {
    $year = [int]$ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 5).Value = $values[$year]
}

# Append six new rows for years 2011-2016.
$row = 93
foreach ($year in 2011..2016) {
    $ws.Cells.Item($row, 1).Value = 320
    $ws.Cells.Item($row, 2).Value = "Guatemala"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $year
    $ws.Cells.Item($row, 5).Value = $values[$year]
    $row++
}
